# Scheduled market-data refresh: update cached price/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 414.57144
$ws.Range("I9").Value = 436.83334
$ws.Range("J9").Value = 397.875
$ws.Range("K9").Value = 436.83334
$ws.Range("L9").Value = 397.875
$ws.Range("M9").Value = -267.83334
$ws.Range("N9").Value = -735.875

$ws.Range("H64").Value = 6557.5
$ws.Range("I64").Value = 5875
$ws.Range("K64").Value = 5875
$ws.Range("M64").Value = -5627

$ws.Range("H67").Value = 6557.5
$ws.Range("I67").Value = 5875
$ws.Range("K67").Value = 5875
$ws.Range("M67").Value = -5017

$ws.Range("H116").Value = 6532.76
$ws.Range("I116").Value = 6333.263
$ws.Range("K116").Value = 6333.263
$ws.Range("M116").Value = -2891.263

$ws.Range("H131").Value = 4267.923
$ws.Range("I131").Value = 2880.6667
$ws.Range("K131").Value = 8642.000100000001
$ws.Range("M131").Value = -3602.000100000001

$ws.Range("H132").Value = 1891.2433
$ws.Range("I132").Value = 1858.6471
$ws.Range("J132").Value = 2260.6667
$ws.Range("K132").Value = 5575.9413
$ws.Range("L132").Value = 6782.000100000001
$ws.Range("M132").Value = -3045.9413
$ws.Range("N132").Value = -11842.0001

$ws.Range("H136").Value = 113771.43
$ws.Range("J136").Value = 113771.43
$ws.Range("L136").Value = 113771.43
$ws.Range("N136").Value = -123971.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41381140
$ws.Range("I32").Value = 42859004
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 42859004
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -42858717
$ws.Range("N32").Value = -1574

$ws.Range("H108").Value = 93167.39999999999
$ws.Range("J108").Value = 93167.39999999999
$ws.Range("L108").Value = 93167.39999999999
$ws.Range("N108").Value = -100847.4

$ws.Range("H111").Value = 97833
$ws.Range("J111").Value = 97833
$ws.Range("L111").Value = 97833
$ws.Range("N111").Value = -106013

$ws.Range("H131").Value = 126000
$ws.Range("J131").Value = 126000
$ws.Range("L131").Value = 126000
$ws.Range("N131").Value = -136080

$ws.Range("H132").Value = 2846.1091
$ws.Range("I132").Value = 2815.4468
$ws.Range("K132").Value = 8446.340400000001
$ws.Range("M132").Value = -5916.340400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41804.23
$ws.Range("I20").Value = 58106.223
$ws.Range("K20").Value = 58106.223
$ws.Range("M20").Value = -57859.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3483.7778
$ws.Range("I31").Value = 1560.4286
$ws.Range("J31").Value = 3948.0344
$ws.Range("K31").Value = 1560.4286
$ws.Range("L31").Value = 3948.0344
$ws.Range("M31").Value = -1265.4286
$ws.Range("N31").Value = -4538.0344

$ws.Range("H34").Value = 3483.7778
$ws.Range("I34").Value = 1560.4286
$ws.Range("J34").Value = 3948.0344
$ws.Range("K34").Value = 1560.4286
$ws.Range("L34").Value = 3948.0344
$ws.Range("M34").Value = -1358.4286
$ws.Range("N34").Value = -4352.0344

$ws.Range("H58").Value = 2371.2546
$ws.Range("I58").Value = 2240.7827
$ws.Range("J58").Value = 3038.111
$ws.Range("K58").Value = 2240.7827
$ws.Range("L58").Value = 3038.111
$ws.Range("M58").Value = -2037.7827
$ws.Range("N58").Value = -3444.111

$ws.Range("H86").Value = 40660.715
$ws.Range("I86").Value = 31498
$ws.Range("J86").Value = 42187.832
$ws.Range("K86").Value = 31498
$ws.Range("L86").Value = 42187.832
$ws.Range("M86").Value = -30375
$ws.Range("N86").Value = -44433.832

$ws.Range("H89").Value = 40660.715
$ws.Range("I89").Value = 31498
$ws.Range("J89").Value = 42187.832
$ws.Range("K89").Value = 157490
$ws.Range("L89").Value = 210939.16
$ws.Range("M89").Value = -151874
$ws.Range("N89").Value = -222171.16

$ws.Range("H99").Value = 2699.6
$ws.Range("I99").Value = 2624.5
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2624.5
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1126.5
$ws.Range("N99").Value = -5996

$ws.Range("H110").Value = 111999
$ws.Range("J110").Value = 111999
$ws.Range("L110").Value = 111999
$ws.Range("N110").Value = -120179

$ws.Range("H126").Value = 2699.6
$ws.Range("I126").Value = 2624.5
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7873.5
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5403.5
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 3931.2856
$ws.Range("I132").Value = 3900.5
$ws.Range("J132").Value = 4116
$ws.Range("K132").Value = 11701.5
$ws.Range("L132").Value = 12348
$ws.Range("M132").Value = -9171.5
$ws.Range("N132").Value = -17408

$ws.Range("H136").Value = 2371.2546
$ws.Range("I136").Value = 2240.7827
$ws.Range("J136").Value = 3038.111
$ws.Range("K136").Value = 6722.348100000001
$ws.Range("L136").Value = 9114.332999999999
$ws.Range("M136").Value = -4172.348100000001
$ws.Range("N136").Value = -14214.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 3298.0667
$ws.Range("I18").Value = 2347.2
$ws.Range("K18").Value = 7041.599999999999
$ws.Range("M18").Value = -6872.599999999999

$ws.Range("H39").Value = 3583.6
$ws.Range("J39").Value = 3583.6
$ws.Range("L39").Value = 10750.8
$ws.Range("N39").Value = -11338.8

$ws.Range("H62").Value = 4917.1113
$ws.Range("J62").Value = 4909.5
$ws.Range("L62").Value = 14728.5
$ws.Range("N62").Value = -16100.5

$ws.Range("H65").Value = 4917.1113
$ws.Range("J65").Value = 4909.5
$ws.Range("L65").Value = 44185.5
$ws.Range("N65").Value = -51049.5

$ws.Range("H121").Value = 9000778
$ws.Range("I121").Value = 617.8
$ws.Range("J121").Value = 18000938
$ws.Range("K121").Value = 1853.4
$ws.Range("L121").Value = 54002814
$ws.Range("M121").Value = -543.3999999999999
$ws.Range("N121").Value = -54005434

$ws.Range("H122").Value = 412
$ws.Range("I122").Value = 474
$ws.Range("J122").Value = 350
$ws.Range("K122").Value = 4266
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -1816
$ws.Range("N122").Value = -8050

$ws.Range("H134").Value = 5105.222
$ws.Range("I134").Value = 5105.222
$ws.Range("K134").Value = 15315.666
$ws.Range("M134").Value = -10245.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 114000
$ws.Range("J119").Value = 114000
$ws.Range("L119").Value = 114000
$ws.Range("N119").Value = -123676

$ws.Range("H126").Value = 3580.1667
$ws.Range("I126").Value = 3695
$ws.Range("J126").Value = 3522.75
$ws.Range("K126").Value = 11085
$ws.Range("L126").Value = 10568.25
$ws.Range("M126").Value = -8615
$ws.Range("N126").Value = -15508.25

$ws.Range("H138").Value = 89915.664
$ws.Range("J138").Value = 89915.664
$ws.Range("L138").Value = 89915.664
$ws.Range("N138").Value = -100195.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6027
$ws.Range("I7").Value = 5867.5
$ws.Range("K7").Value = 5867.5
$ws.Range("M7").Value = -5755.5

$ws.Range("H40").Value = 4600.1113
$ws.Range("I40").Value = 3932.2856
$ws.Range("J40").Value = 6937.5
$ws.Range("K40").Value = 3932.2856
$ws.Range("L40").Value = 6937.5
$ws.Range("M40").Value = -3796.2856
$ws.Range("N40").Value = -7209.5

$ws.Range("H46").Value = 7217.5625
$ws.Range("I46").Value = 3479.6
$ws.Range("J46").Value = 7909.778
$ws.Range("K46").Value = 3479.6
$ws.Range("L46").Value = 7909.778
$ws.Range("M46").Value = -3291.6
$ws.Range("N46").Value = -8285.778

$ws.Range("H116").Value = 158999.5
$ws.Range("J116").Value = 158999.5
$ws.Range("L116").Value = 158999.5
$ws.Range("N116").Value = -168177.5

$ws.Range("H122").Value = 22805
$ws.Range("I122").Value = 26491.834
$ws.Range("K122").Value = 79475.50199999999
$ws.Range("M122").Value = -77025.50199999999

$ws.Range("H126").Value = 6027
$ws.Range("I126").Value = 5867.5
$ws.Range("K126").Value = 17602.5
$ws.Range("M126").Value = -15132.5

$ws.Range("H132").Value = 5029.6216
$ws.Range("I132").Value = 4656.826
$ws.Range("K132").Value = 13970.478
$ws.Range("M132").Value = -11440.478

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2607.389
$ws.Range("I132").Value = 2475.2666
$ws.Range("J132").Value = 3268
$ws.Range("K132").Value = 7425.7998
$ws.Range("L132").Value = 9804
$ws.Range("M132").Value = -4895.7998
$ws.Range("N132").Value = -14864
